$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.235.95"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.827.75"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'235.41"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").Value = "'0.6002"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.06956"
$ws.Range("E8").Value = "  -2.11%  "
$ws.Range("D9").Value = "'0.2768"
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("D10").Value = "'23.49"
$ws.Range("E10").Value = "  -2.82%  "
$ws.Range("D11").Value = "'0.07612"
$ws.Range("E11").Value = "  -1.07%  "
$ws.Range("D12").Value = "1.834.90"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").Value = "'4.739"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").Value = "'0.6312"
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("D15").Value = "'0.000009845"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").Value = "'77.58"
$ws.Range("E16").Value = "  -2.33%  "
$ws.Range("D17").Value = "28.979.64"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").Value = "'5.559"
$ws.Range("E18").Value = "  -8.38%  "
$ws.Range("D19").Value = "'216.54"
$ws.Range("E19").Value = "  -5.85%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "'11.57"
$ws.Range("E21").Value = "  -1.70%  "
$ws.Range("D22").Value = "'6.874"
$ws.Range("E22").Value = "  -1.99%  "
$ws.Range("D23").Value = "'1.006"
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").Value = "'156.25"
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").Value = "'7.956"
$ws.Range("E25").Value = "  -1.87%  "
$ws.Range("D26").Value = "'0.1287"
$ws.Range("E26").Value = "  +1.36%  "
$ws.Range("D27").Value = "'16.52"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").Value = "'0.06435"
$ws.Range("E28").Value = "  -5.30%  "
$ws.Range("D29").Value = "'1.413"
$ws.Range("E29").Value = "  -4.07%  "
$ws.Range("D30").Value = "'1.443"
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("D31").Value = "'3.827"
$ws.Range("E31").Value = "  +1.69%  "
$ws.Range("D32").Value = "'3.787"
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").Value = "'1.095"
$ws.Range("E33").Value = "  -3.50%  "
$ws.Range("D34").Value = "'1.724"
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("D35").Value = "'0.6474"
$ws.Range("E35").Value = "  -2.20%  "
$ws.Range("D36").Value = "'2.545"
$ws.Range("E36").Value = "  -1.08%  "
$ws.Range("D37").Value = "'2.758"
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("E38").Value = "  -0.56%  "
$ws.Range("D39").Value = "'6.608"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "1.137.77"
$ws.Range("E40").Value = "  -6.84%  "
$ws.Range("D41").Value = "'0.8931"
$ws.Range("E41").Value = "  -3.70%  "
$ws.Range("D42").Value = "'1.004"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").Value = "1.999.75"
$ws.Range("E43").Value = "  +0.67%  "
$ws.Range("D44").Value = "'100.61"
$ws.Range("E44").Value = "  +1.50%  "
$ws.Range("D45").Value = "'62.17"
$ws.Range("E45").Value = "  -2.04%  "
$ws.Range("D46").Value = "'0.00000000113"
$ws.Range("E46").Value = "  -4.11%  "
$ws.Range("D47").Value = "'1.620"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("D48").Value = "'8.502"
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("D49").Value = "'0.05501"
$ws.Range("E49").Value = "  -1.89%  "
$ws.Range("D50").Value = "'0.4538"
$ws.Range("E50").Value = "  -0.91%  "
$ws.Range("D51").Value = "'6.384"
$ws.Range("E51").Value = "  -3.05%  "
